# Add a new "before" worksheet as the last tab, make it the active/selected
# sheet (matching activeTab="4" / tabSelected="1" moving off the first
# sheet), and populate it with the new JETT tag-listener "before callback"
# sample content that was added to sharedStrings.xml.

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last worksheet so it lands at the
# end of the tab strip (sheetId=5, rId5, appears after "implInstance").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "before"

# Populate cells in an order that reproduces the shared-string table order
# from the target workbook: "Number of Employees:" (A1) first, then the
# "before" replacement text (B2), then the "will be replaced" text (B1).
$newSheet.Range("A1").Value = "Number of Employees:"
$newSheet.Range("B2").Value = '<jt:for var="x" start="1" end="1" onProcessed="${boldTagListener}">The above will be replaced by ${employees.size()}</jt:for>'
$newSheet.Range("B1").Value = '<jt:for var="x" start="1" end="1" onProcessed="${boldTagListener}">Will be replaced!</jt:for>'

# Column A should be wide enough to fully show "Number of Employees:"
# (mirrors the bestFit/customWidth column sizing on the source sheet).
$newSheet.Columns.Item(1).ColumnWidth = 20.833333333333332

# Make the new "before" sheet the active/selected tab, like the diff shows
# (tabSelected moves from the first sheet to this new last sheet, and the
# workbook's activeTab points at index 4).
$newSheet.Activate()
